# data run for sg_rr_100_030 2023-12-08 16-08-32
#
# Adds the next row of analysis results (row 35) to the FSR data-analysis
# table on Sheet1, mirroring the existing rows' layout:
#   A: Data CSV filename
#   B: Wavelength step size/nm
#   C: Start array index
#   D: End array index
#   E: Start wavelength/nm
#   F: End wavelength/nm
#   G: prominence/dBm
#   H: distance (note text)
#   I: approx_fsr/nm
#   J: fsr_mean/nm
#   K: fsr_std error/nm
#   L: double count check passed?
#   M: free-text note about how prominence was chosen for this run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

$ws.Cells.Item($row, 1).Value  = "sg_rr_100_030 2023-12-08 16-08-32.csv"
$ws.Cells.Item($row, 2).Value  = 0.01
$ws.Cells.Item($row, 3).Value  = 1000
$ws.Cells.Item($row, 4).Value  = 5001
$ws.Cells.Item($row, 5).Value  = 1530
$ws.Cells.Item($row, 6).Value  = 1570
$ws.Cells.Item($row, 7).Value  = 0.5
$ws.Cells.Item($row, 8).Value  = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = 0.98128205128205004
$ws.Cells.Item($row, 11).Value = 0.00316397329552258
$ws.Cells.Item($row, 12).Value = "yes"
$ws.Cells.Item($row, 13).Value = "prominence set by looking at roughly biggest height span of noise bits that don't appear visually to contain resonance peaks"

# Match the saved selection/scroll position from the commit (K35 active,
# scrolled so row 26 / column B is the top-left visible cell).
$ws.Range("K35").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 2
